$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Find the last used row (data rows 2..47, header row 1)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Add the three new header cells (Wins, Losses, Ties) in AD1:AF1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style from an existing header cell (A1) to the new headers
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Re-set the values after the style paste (PasteSpecial formats only shouldn't
# touch values, but set them again defensively)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the team record (Wins=96, Losses=66, Ties=0) for every data row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 96
    $ws.Cells.Item($r, 31).Value = 66
    $ws.Cells.Item($r, 32).Value = 0
}
